# Add the three new character styles used by the new paragraphs.
$d = $word.ActiveDocument

$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# Apply GaNStyle to every "Kampagnendaten 2022 ..." heading run (4 occurrences).
$rng = $d.Content
while ($rng.Find.Execute("Kampagnendaten 2022 für das Sternbild Löwe: 14.-23. April, 14.-23. Mai", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
}

# Apply GaNParagraph to the intro paragraph run.
$rng2 = $d.Content
if ($rng2.Find.Execute("Mach mit an einer weltweiten Kampagne, die schwächsten sichtbaren Sterne zu beobachten und aufzuzeichnen, um die Lichtverschmutzung an einem Ort zu messen. Durch das Auffinden und Beobachten des Sternbild Löwe am Nachthimmel und den Vergleich mit den Helligkeitskarten, lernen Menschen auf der ganzen Erde, wie die Lichter in ihrer Gemeinde zur Lichtverschmutzung beitragen. Dein Beitrag zur Online-Datenbank beschreibt den sichtbaren Nachthimmel.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng2.Style = "GaNParagraph"
}

# Apply GaNLinks to the credit/link paragraph run.
$rng3 = $d.Content
if ($rng3.Find.Execute("Die Schaubilder in diesem Dokument wurden von Jan Hollan, CzechGlobe, bereitgestellt. (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng3.Style = "GaNLinks"
}
